$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2026-01-26 Monday" "2026-01-27 Tuesday"

Replace-Text "344×7=2408" "367×9=3303"
Replace-Text "686×6=4116" "422×3=1266"
Replace-Text "986×5=4930" "223×3=669"
Replace-Text "193×2=386" "800×4=3200"
Replace-Text "996×2=1992" "571×5=2855"

Replace-Text "576×2=1152" "183×9=1647"
Replace-Text "405×8=3240" "568×4=2272"
Replace-Text "463×8=3704" "666×6=3996"
Replace-Text "998×8=7984" "293×9=2637"
Replace-Text "460×5=2300" "944×4=3776"

Replace-Text "934×6=5604" "281×5=1405"
Replace-Text "199×7=1393" "988×2=1976"
Replace-Text "613×8=4904" "623×7=4361"
Replace-Text "881×8=7048" "903×4=3612"
Replace-Text "667×2=1334" "294×4=1176"

Replace-Text "391×3=1173" "927×4=3708"
Replace-Text "238×2=476" "523×4=2092"
Replace-Text "803×8=6424" "840×4=3360"
Replace-Text "333×4=1332" "310×4=1240"
Replace-Text "967×2=1934" "415×9=3735"

Replace-Text "131×6=786" "288×2=576"
Replace-Text "442×7=3094" "467×3=1401"
Replace-Text "379×8=3032" "268×4=1072"
Replace-Text "993×2=1986" "265×3=795"
Replace-Text "208×9=1872" "545×9=4905"
